# Update cryptos list with latest prices / 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '43.496.67'
$ws.Range("E2").Value = '  -0.78%  '
# Row 3 - Ethereum
$ws.Range("D3").Value = '2.275.36'
$ws.Range("E3").Value = '  -0.88%  '
# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.36%  '
# Row 5 - Solana
$ws.Range("D5").Value = "'122.59"
$ws.Range("E5").Value = '  +6.33%  '
# Row 6 - BNB
$ws.Range("D6").Value = "'265.61"
$ws.Range("E6").Value = '  -1.35%  '
# Row 7 - XRP
$ws.Range("E7").Value = '  +1.83%  '
# Row 8 - USDC
$ws.Range("E8").Value = '  +0.23%  '
# Row 9 - Cardano
$ws.Range("E9").Value = '  +0.89%  '
# Row 10 - Avalanche
$ws.Range("D10").Value = "'48.19"
$ws.Range("E10").Value = '  -0.73%  '
# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0945"
$ws.Range("E11").Value = '  +0.27%  '
# Row 12 - Polkadot
$ws.Range("D12").Value = "'9.03"
$ws.Range("E12").Value = '  -0.21%  '
# Row 13 - TRON
$ws.Range("E13").Value = '  -0.70%  '
# Row 14 - Chainlink
$ws.Range("D14").Value = "'15.39"
$ws.Range("E14").Value = '  -2.77%  '
# Row 15 - Polygon
$ws.Range("D15").Value = "'0.891"
$ws.Range("E15").Value = '  +4.06%  '
# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '2.619.76'
$ws.Range("E16").Value = '  -0.78%  '
# Row 17 - WrappedEther
$ws.Range("D17").Value = '2.279.50'
$ws.Range("E17").Value = '  -1.01%  '
# Row 18 - WrappedBTC
$ws.Range("D18").Value = '43.634.55'
$ws.Range("E18").Value = '  -0.23%  '
# Row 19 - ShibaInu
$ws.Range("E19").Value = '  -0.20%  '
# Row 20 - Uniswap
$ws.Range("D20").Value = "'6.99"
$ws.Range("E20").Value = '  -1.98%  '
# Row 21 - Litecoin
$ws.Range("D21").Value = "'72.24"
$ws.Range("E21").Value = '  -0.19%  '
# Row 22 - ImmutableX
$ws.Range("E22").Value = '  +0.18%  '
# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'235.79"
$ws.Range("E23").Value = '  +1.13%  '
# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'9.52"
$ws.Range("E24").Value = '  -3.33%  '
# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = '  -3.13%  '
# Row 26 - Dai
$ws.Range("E26").Value = '  +1.83%  '
# Row 27 - Cosmos
$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = '  +1.00%  '
# Row 28 - InjectiveProtocol
$ws.Range("D28").Value = "'42.22"
$ws.Range("E28").Value = '  -2.07%  '
# Row 29 - WEMIXToken
$ws.Range("E29").Value = '  -0.62%  '
# Row 30 - Toncoin
$ws.Range("E30").Value = '  -0.09%  '
# Row 31 - Monero
$ws.Range("D31").Value = "'171.70"
$ws.Range("E31").Value = '  -1.97%  '
# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'21.70"
$ws.Range("E32").Value = '  +0.32%  '
# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0906"
$ws.Range("E33").Value = '  -3.21%  '
# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = '  +0.53%  '
# Row 35 - Stellar
$ws.Range("E35").Value = '  +1.97%  '
# Row 36 - VeChain
$ws.Range("D36").Value = "'0.0378"
$ws.Range("E36").Value = '  +4.43%  '
# Row 37 - RenderToken
$ws.Range("D37").Value = "'4.64"
$ws.Range("E37").Value = '  -3.77%  '
# Row 38 - NEARProtocol
$ws.Range("D38").Value = "'4.03"
$ws.Range("E38").Value = '  +5.07%  '
# Row 39 - Kaspa
$ws.Range("D39").Value = "'0.108"
$ws.Range("E39").Value = '  +1.54%  '
# Row 40 - LidoDAOToken
$ws.Range("E40").Value = '  +5.16%  '
# Row 41 - was MultiversX, now Celestia (rows 41/42 swapped order)
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = "'14.07"
$ws.Range("E41").Value = '  -1.95%  '
# Row 42 - was Celestia, now MultiversX
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = "'74.97"
$ws.Range("E42").Value = '  -1.27%  '
# Row 43 - Algorand
$ws.Range("E43").Value = '  -2.01%  '
# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = '  -0.03%  '
# Row 45 - ARBITRUM
$ws.Range("E45").Value = '  -3.19%  '
# Row 46 - THORChain
$ws.Range("E46").Value = '  -9.46%  '
# Row 47 - was TrustWalletToken, now ordi (rows 47/48 swapped order)
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").Value = "'74.04"
$ws.Range("E47").Value = '  +37.19%  '
# Row 48 - was ordi, now TrustWalletToken
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = "'1.27"
$ws.Range("E48").Value = '  +0.75%  '
# Row 49 - FraxShare
$ws.Range("E49").Value = '  -3.31%  '
# Row 50 - Cronos
$ws.Range("E50").Value = '  +0.93%  '
# Row 51 - Aave
$ws.Range("D51").Value = "'101.69"
$ws.Range("E51").Value = '  -0.51%  '

